$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    # Numeric / date-looking strings need a leading apostrophe so Excel
    # stores them as text instead of coercing to a number/date; then strip
    # the formatting it stamps on for the quote-prefix so no style index
    # is left behind on the cell.
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# --- Row 20: the two previously-blank cells (B20, D20) are removed ---
$ws.Cells.Item(20, 2).ClearContents()
$ws.Cells.Item(20, 4).ClearContents()

# --- New rows 21-24 ---
$rows = @(
    @{ A=20; B="david"; C="2020-01-02"; D="pmma"; E="Cut"; F="90"; G="90"; H="900"; I="5000"; J="1/0"; K="1"; L="Enter here useful comments for the future" },
    @{ A=21; B="david"; C="2020-01-02"; D="pmma"; E="Cut"; F="90"; G="90"; H="900"; I="5000"; J="1/0"; K="1"; L="Enter here useful comments for the future" },
    @{ A=22; B="sads";  C="2020-01-02"; D="sd";   E="Cut"; F="90"; G="90"; H="900"; I="5000"; J="1/0"; K="1"; L="Enter here useful comments for the future" },
    @{ A=23; B="aS";    C="2020-01-02"; D="Asa";  E="Cut"; F="90"; G="90"; H="900"; I="5000"; J="1/0"; K="1"; L="Enter here useful comments for the future" }
)

$r = 21
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    Set-TextCell $r 2  $row.B
    Set-TextCell $r 3  $row.C
    Set-TextCell $r 4  $row.D
    Set-TextCell $r 5  $row.E
    Set-TextCell $r 6  $row.F
    Set-TextCell $r 7  $row.G
    Set-TextCell $r 8  $row.H
    Set-TextCell $r 9  $row.I
    Set-TextCell $r 10 $row.J
    Set-TextCell $r 11 $row.K
    Set-TextCell $r 12 $row.L
    $r++
}
